$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-10-13 Sunday" "2024-10-14 Monday"

Replace-Text "930×7=" "494×2="
Replace-Text "608×2=" "800×9="
Replace-Text "831×2=" "267×3="
Replace-Text "347×8=" "186×8="
Replace-Text "274×2=" "664×4="

Replace-Text "113×6=" "721×9="
Replace-Text "930×8=" "770×2="
Replace-Text "462×6=" "199×4="
Replace-Text "545×8=" "554×8="
Replace-Text "384×3=" "542×2="

Replace-Text "103×4=" "616×4="
Replace-Text "408×7=" "442×4="
Replace-Text "432×6=" "590×4="
Replace-Text "675×9=" "523×7="
Replace-Text "665×3=" "960×2="

Replace-Text "840×9=" "393×2="
Replace-Text "247×2=" "693×8="
Replace-Text "461×2=" "714×5="
Replace-Text "747×5=" "582×6="
Replace-Text "877×2=" "156×8="

Replace-Text "957×7=" "950×5="
Replace-Text "159×3=" "950×4="
Replace-Text "603×5=" "190×8="
Replace-Text "588×5=" "849×7="
Replace-Text "712×8=" "147×8="
